# Regenerate the lattice-multiplication practice grid with a new set of
# problems, keeping the existing 5x3 table (and per-cell layout: product
# header, factor-digit line, dash rule, two partial-product stub lines)
# untouched structurally -- only the numbers change.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each entry: row, col, header "A x B", digit line, two stub lines.
# The dashed separator line ("  ----") never changes, so it's filled in
# below rather than repeated in the table.
$NL = [char]11   # soft line break (w:br), matches Range.Text's run-break char

$cells = @(
    @{R=1; C=1; H="29 x 94"; D="  9    4"; L1="2|    |"; L2="9|    |"},
    @{R=1; C=2; H="22 x 59"; D="  5    9"; L1="2|    |"; L2="2|    |"},
    @{R=1; C=3; H="73 x 13"; D="  1    3"; L1="7|    |"; L2="3|    |"},

    @{R=2; C=1; H="94 x 19"; D="  1    9"; L1="9|    |"; L2="4|    |"},
    @{R=2; C=2; H="64 x 43"; D="  4    3"; L1="6|    |"; L2="4|    |"},
    @{R=2; C=3; H="71 x 69"; D="  6    9"; L1="7|    |"; L2="1|    |"},

    @{R=3; C=1; H="41 x 78"; D="  7    8"; L1="4|    |"; L2="1|    |"},
    @{R=3; C=2; H="27 x 19"; D="  1    9"; L1="2|    |"; L2="7|    |"},
    @{R=3; C=3; H="83 x 97"; D="  9    7"; L1="8|    |"; L2="3|    |"},

    @{R=4; C=1; H="42 x 44"; D="  4    4"; L1="4|    |"; L2="2|    |"},
    @{R=4; C=2; H="50 x 35"; D="  3    5"; L1="5|    |"; L2="0|    |"},
    @{R=4; C=3; H="34 x 72"; D="  7    2"; L1="3|    |"; L2="4|    |"},

    @{R=5; C=1; H="86 x 98"; D="  9    8"; L1="8|    |"; L2="6|    |"},
    @{R=5; C=2; H="45 x 23"; D="  2    3"; L1="4|    |"; L2="5|    |"},
    @{R=5; C=3; H="69 x 87"; D="  8    7"; L1="6|    |"; L2="9|    |"}
)

foreach ($c in $cells) {
    $newText = $c.H + $NL + $c.D + $NL + "  ----" + $NL + $c.L1 + $NL + $c.L2
    $t.Cell($c.R, $c.C).Range.Text = $newText
}

Write-Output "Lattice multiplication exercises updated."
